# Apply the edits described by the diff using Find/Replace across the whole document.
$d = $word.ActiveDocument

# Each pair is: old text, new text
$replacements = @(
    @("2024-08-14 Wednesday", "2024-08-15 Thursday"),
    @("959÷9=106, 5", "132÷5=26, 2"),
    @("955÷2=477, 1", "247÷8=30, 7"),
    @("438÷2=219, 0", "693÷9=77, 0"),
    @("333÷4=83, 1", "868÷6=144, 4"),
    @("612÷9=68, 0", "244÷6=40, 4"),
    @("907÷4=226, 3", "666÷4=166, 2"),
    @("950÷2=475, 0", "507÷6=84, 3"),
    @("440÷3=146, 2", "772÷9=85, 7"),
    @("722÷8=90, 2", "515÷8=64, 3"),
    @("603÷7=86, 1", "705÷5=141, 0"),
    @("927÷5=185, 2", "123÷8=15, 3"),
    @("583÷2=291, 1", "491÷3=163, 2"),
    @("703÷4=175, 3", "690÷4=172, 2"),
    @("236÷8=29, 4", "468÷5=93, 3"),
    @("988÷8=123, 4", "629÷8=78, 5"),
    @("265÷5=53, 0", "423÷5=84, 3"),
    @("760÷9=84, 4", "695÷2=347, 1"),
    @("421÷7=60, 1", "193÷3=64, 1"),
    @("755÷7=107, 6", "420÷7=60, 0"),
    @("218÷9=24, 2", "328÷3=109, 1"),
    @("387÷9=43, 0", "387÷3=129, 0"),
    @("483÷9=53, 6", "153÷8=19, 1"),
    @("927÷2=463, 1", "472÷8=59, 0"),
    @("280÷8=35, 0", "549÷4=137, 1"),
    @("317÷6=52, 5", "439÷7=62, 5")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
